$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E data cells keep their original text type (avoid Excel
# auto-converting numeric-looking strings like "306.60" into numbers,
# which would drop the trailing zero on save).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.225.01"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.270.14"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "306.60"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "97.58"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "35.44"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "6.87"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "2.620.24"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "14.72"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "2.276.46"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "42.102.98"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "12.33"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "67.77"
$ws.Range("D23").Value = "237.20"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "23.55"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "37.42"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("D29").Value = "9.60"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "162.68"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "5.27"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D35").Value = "17.73"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").Value = "0.0736"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "4.11"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "19.03"
$ws.Range("E44").Value = "  -2.64%  "
$ws.Range("D45").Value = "0.0281"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").Value = "9.97"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "53.74"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "2.491.40"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "72.20"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "92.53"
$ws.Range("E51").Value = "  +0.14%  "
